$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 8413.857
$ws.Range("I4").Value = 8413.857
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 8413.857
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -8299.857

$ws.Range("H32").Value = 19000
$ws.Range("I32").Value = 19000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 19000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -18674

$ws.Range("H41").Value = 1316.0588
$ws.Range("I41").Value = 1375.6666
$ws.Range("J41").Value = 1249
$ws.Range("K41").Value = 1375.6666
$ws.Range("L41").Value = 1249
$ws.Range("M41").Value = -935.6666
$ws.Range("N41").Value = -2129

$ws.Range("H55").Value = 194.33333
$ws.Range("I55").Value = 205.88889
$ws.Range("J55").Value = 159.66667
$ws.Range("K55").Value = 205.88889
$ws.Range("L55").Value = 159.66667
$ws.Range("M55").Value = 8.111109999999996
$ws.Range("N55").Value = -587.6666700000001

$ws.Range("H58").Value = 935.5
$ws.Range("I58").Value = 704.3333
$ws.Range("J58").Value = 1166.6666
$ws.Range("K58").Value = 2112.9999
$ws.Range("L58").Value = 3499.9998
$ws.Range("M58").Value = -1962.9999
$ws.Range("N58").Value = -3799.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 7000
$ws.Range("I12").Value = 5500
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 5500
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -5327
$ws.Range("N12").Value = -10346

$ws.Range("H63").Value = 4712.5
$ws.Range("I63").Value = 4232.3335
$ws.Range("J63").Value = 6153
$ws.Range("K63").Value = 4232.3335
$ws.Range("L63").Value = 6153
$ws.Range("M63").Value = -3546.3335

$ws.Range("H66").Value = 4712.5
$ws.Range("I66").Value = 4232.3335
$ws.Range("J66").Value = 6153
$ws.Range("K66").Value = 21161.6675
$ws.Range("L66").Value = 30765
$ws.Range("M66").Value = -17729.6675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 2750
$ws.Range("I46").Value = 2750
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2750
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2452

$ws.Range("H63").Value = 119499.5
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 119499.5
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 119499.5
$ws.Range("N63").Value = -120871.5

$ws.Range("H66").Value = 119499.5
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 119499.5
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 358498.5
$ws.Range("N66").Value = -365362.5

$ws.Range("H82").Value = 27121.666
$ws.Range("I82").Value = 4250
$ws.Range("J82").Value = 49993.332
$ws.Range("K82").Value = 4250
$ws.Range("L82").Value = 49993.332
$ws.Range("M82").Value = -3867
$ws.Range("N82").Value = -50759.332

$ws.Range("H85").Value = 27121.666
$ws.Range("I85").Value = 4250
$ws.Range("J85").Value = 49993.332
$ws.Range("K85").Value = 4250
$ws.Range("L85").Value = 49993.332
$ws.Range("M85").Value = -2924
$ws.Range("N85").Value = -52645.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2670.5386
$ws.Range("I7").Value = 90.2
$ws.Range("J7").Value = 4283.25
$ws.Range("K7").Value = 90.2
$ws.Range("L7").Value = 4283.25
$ws.Range("M7").Value = 22.8

$ws.Range("H10").Value = 2575
$ws.Range("I10").Value = 4000
$ws.Range("J10").Value = 2100
$ws.Range("K10").Value = 4000
$ws.Range("L10").Value = 2100
$ws.Range("M10").Value = -3861
$ws.Range("N10").Value = -2378

$ws.Range("H22").Value = 666.6667
$ws.Range("I22").Value = 666.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 666.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -316.6667
$ws.Range("N22").ClearContents()

$ws.Range("H105").Value = 2186.1785
$ws.Range("I105").Value = 1265.7059
$ws.Range("J105").Value = 3608.7273
$ws.Range("K105").Value = 1265.7059
$ws.Range("L105").Value = 3608.7273
$ws.Range("M105").Value = 481.2941000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 116.588234
$ws.Range("I2").Value = 65.5
$ws.Range("J2").Value = 144.45454
$ws.Range("K2").Value = 393
$ws.Range("L2").Value = 866.7272400000001
$ws.Range("M2").Value = -280
$ws.Range("N2").Value = -1092.72724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 592.8333
$ws.Range("I2").Value = 175.11111
$ws.Range("J2").Value = 1846
$ws.Range("K2").Value = 175.11111
$ws.Range("L2").Value = 1846
$ws.Range("M2").Value = -62.11111
$ws.Range("N2").Value = -2072

$ws.Range("H24").Value = 2875
$ws.Range("I24").Value = 2500
$ws.Range("J24").Value = 4000
$ws.Range("K24").Value = 2500
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = -2327
$ws.Range("N24").Value = -4346

$ws.Range("H26").Value = 9900
$ws.Range("I26").Value = 9750
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 9750
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = -9470
$ws.Range("N26").Value = -10560

$ws.Range("H50").Value = 9900
$ws.Range("I50").Value = 9750
$ws.Range("J50").Value = 10000
$ws.Range("K50").Value = 9750
$ws.Range("L50").Value = 10000
$ws.Range("M50").Value = -9252
$ws.Range("N50").Value = -10996

$ws.Range("H59").Value = 3333.3333
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 3333.3333
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 3333.3333
$ws.Range("N59").Value = -4499.3333

$ws.Range("H122").Value = 2067.3
$ws.Range("I122").Value = 2019.7778
$ws.Range("J122").Value = 2495
$ws.Range("K122").Value = 6059.3334
$ws.Range("L122").Value = 7485
$ws.Range("M122").Value = -3609.3334
$ws.Range("N122").Value = -12385

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H19").Value = 250
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 250
$ws.Range("N19").Value = -598

$ws.Range("H41").Value = 21657.334
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 21657.334
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 21657.334
$ws.Range("N41").Value = -22437.334

$ws.Range("H61").Value = 20519.666
$ws.Range("I61").Value = 10224.5
$ws.Range("J61").Value = 41110
$ws.Range("K61").Value = 10224.5
$ws.Range("L61").Value = 41110
$ws.Range("M61").Value = -9932.5
$ws.Range("N61").Value = -41694

$ws.Range("H139").Value = 75333.336
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 75333.336
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 75333.336
$ws.Range("N139").Value = -85613.336
